# "Comprobación de los posibles errores de lectura"
#
# Changes applied (per the target diff):
#  1. Workbook window height metadata (bookViews/workbookView/@windowHeight
#     3720 -> 4650) - cosmetic app-window size, attempted below for
#     completeness even though it is not part of the worksheet contents.
#  2. Sheet1 ("Hoja1") selection/active cell moves from D4 to D5.
#  3. Cell D4 on Hoja1 gets the literal value -9 (was empty).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# 3. Fill in D4 with -9 (keeps its existing style s="1").
$ws.Range("D4").Value = -9

# 2. Move the selection / active cell to D5.
$ws.Range("D5").Select()

# 1. Best-effort: record the new application window height.
$excel.ActiveWindow.Height = 4650
